$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 41589
$ws.Range("D2").Value = 60048447
$ws.Range("C3").Value = 98772
$ws.Range("D3").Value = 144683350
$ws.Range("C4").Value = 33624
$ws.Range("D4").Value = 49768313
$ws.Range("C5").Value = 9616
$ws.Range("D5").Value = 14284242
$ws.Range("C6").Value = 2379
$ws.Range("D6").Value = 3534471
$ws.Range("C7").Value = 249
$ws.Range("D7").Value = 368593
$ws.Range("C12").Value = 44732
$ws.Range("D12").Value = 60515782
$ws.Range("C13").Value = 10534
$ws.Range("D13").Value = 15212918
$ws.Range("C14").Value = 27819
$ws.Range("D14").Value = 40765897
$ws.Range("C15").Value = 8818
$ws.Range("D15").Value = 13085166
$ws.Range("C16").Value = 2335
$ws.Range("D16").Value = 3467875
$ws.Range("C20").Value = 10976
$ws.Range("D20").Value = 14449504
$ws.Range("C21").Value = 14479
$ws.Range("D21").Value = 20865750
$ws.Range("C22").Value = 33662
$ws.Range("D22").Value = 49355296
$ws.Range("C23").Value = 10835
$ws.Range("D23").Value = 16099239
$ws.Range("C24").Value = 2876
$ws.Range("D24").Value = 4274115
$ws.Range("C26").Value = 44
$ws.Range("D26").Value = 65453
$ws.Range("C27").Value = 12503
$ws.Range("D27").Value = 16615544
$ws.Range("C28").Value = 8423
$ws.Range("D28").Value = 12181684
$ws.Range("C29").Value = 24221
$ws.Range("D29").Value = 35524253
$ws.Range("C30").Value = 8334
$ws.Range("D30").Value = 12389707
$ws.Range("C31").Value = 2127
$ws.Range("D31").Value = 3172708
$ws.Range("C32").Value = 430
$ws.Range("D32").Value = 635621
$ws.Range("C34").Value = 8992
$ws.Range("D34").Value = 11844848
$ws.Range("C35").Value = 3684
$ws.Range("D35").Value = 5314068
$ws.Range("C36").Value = 8631
$ws.Range("D36").Value = 12611045
$ws.Range("C37").Value = 3404
$ws.Range("D37").Value = 5047008
$ws.Range("C38").Value = 875
$ws.Range("D38").Value = 1303555
$ws.Range("C39").Value = 181
$ws.Range("D39").Value = 269186
$ws.Range("C41").Value = 2761
$ws.Range("D41").Value = 3715647
$ws.Range("C42").Value = 18834
$ws.Range("D42").Value = 27191903
$ws.Range("C43").Value = 54807
$ws.Range("D43").Value = 80290965
$ws.Range("C44").Value = 20098
$ws.Range("D44").Value = 29832081
$ws.Range("C46").Value = 1411
$ws.Range("D46").Value = 2106144
$ws.Range("C50").Value = 18270
$ws.Range("D50").Value = 24186994
$ws.Range("C51").Value = 2370
$ws.Range("D51").Value = 3438526
$ws.Range("C52").Value = 7981
$ws.Range("D52").Value = 11723028
$ws.Range("C53").Value = 2675
$ws.Range("D53").Value = 3992633
$ws.Range("C54").Value = 843
$ws.Range("D54").Value = 1259414
$ws.Range("C57").Value = 7953
$ws.Range("D57").Value = 10941837
$ws.Range("C58").Value = 1650
$ws.Range("D58").Value = 3316440
$ws.Range("C59").Value = 3939
$ws.Range("D59").Value = 7882312
$ws.Range("C60").Value = 1557
$ws.Range("D60").Value = 3121294
$ws.Range("C61").Value = 528
$ws.Range("D61").Value = 1053083
$ws.Range("C62").Value = 191
$ws.Range("D62").Value = 394487
$ws.Range("C64").Value = 2561
$ws.Range("D64").Value = 4752864
$ws.Range("C65").Value = 16976
$ws.Range("D65").Value = 24498817
$ws.Range("C66").Value = 48309
$ws.Range("D66").Value = 70617001
$ws.Range("C67").Value = 16882
$ws.Range("D67").Value = 25083962
$ws.Range("C69").Value = 1091
$ws.Range("D69").Value = 1622199
$ws.Range("C73").Value = 16130
$ws.Range("D73").Value = 21164772
$ws.Range("C74").Value = 62563
$ws.Range("D74").Value = 90931790
$ws.Range("C75").Value = 170320
$ws.Range("D75").Value = 250643438
$ws.Range("C76").Value = 72772
$ws.Range("D76").Value = 108376532
$ws.Range("C77").Value = 23812
$ws.Range("D77").Value = 35562074
$ws.Range("C78").Value = 6195
$ws.Range("D78").Value = 9244877
$ws.Range("C79").Value = 451
$ws.Range("D79").Value = 671458
$ws.Range("C85").Value = 61562
$ws.Range("D85").Value = 83003511
$ws.Range("C86").Value = 5133
$ws.Range("D86").Value = 7435873
$ws.Range("C87").Value = 12616
$ws.Range("D87").Value = 18525814
$ws.Range("C88").Value = 4136
$ws.Range("D88").Value = 6163126
$ws.Range("C89").Value = 1455
$ws.Range("D89").Value = 2173611
$ws.Range("C93").Value = 5874
$ws.Range("D93").Value = 7877874
$ws.Range("C94").Value = 1849
$ws.Range("D94").Value = 2661993
$ws.Range("C95").Value = 5865
$ws.Range("D95").Value = 8642757
$ws.Range("C96").Value = 2109
$ws.Range("D96").Value = 3143431
$ws.Range("C98").Value = 223
$ws.Range("D98").Value = 339109
$ws.Range("C99").Value = 25
$ws.Range("D99").Value = 37500
$ws.Range("C101").Value = 3984
$ws.Range("D101").Value = 5280363
$ws.Range("C102").Value = 919
$ws.Range("D102").Value = 1779061
$ws.Range("C103").Value = 615
$ws.Range("D103").Value = 1259284
$ws.Range("C107").Value = 11900
$ws.Range("D107").Value = 17245765
$ws.Range("C108").Value = 31121
$ws.Range("D108").Value = 45675276
$ws.Range("C109").Value = 10429
$ws.Range("D109").Value = 15502760
$ws.Range("C110").Value = 2891
$ws.Range("D110").Value = 4309571
$ws.Range("C115").Value = 10493
$ws.Range("D115").Value = 13803149
$ws.Range("C116").Value = 33414
$ws.Range("D116").Value = 48133995
$ws.Range("C117").Value = 70927
$ws.Range("D117").Value = 103732868
$ws.Range("C118").Value = 22744
$ws.Range("D118").Value = 33782586
$ws.Range("C119").Value = 6526
$ws.Range("D119").Value = 9715052
$ws.Range("C120").Value = 1305
$ws.Range("D120").Value = 1950233
$ws.Range("C121").Value = 127
$ws.Range("D121").Value = 185295
$ws.Range("C125").Value = 27695
$ws.Range("D125").Value = 36870555
$ws.Range("C126").Value = 39919
$ws.Range("D126").Value = 57542589
$ws.Range("C127").Value = 83348
$ws.Range("D127").Value = 121766362
$ws.Range("C128").Value = 25527
$ws.Range("D128").Value = 37877249
$ws.Range("C129").Value = 6967
$ws.Range("D129").Value = 10353009
$ws.Range("C130").Value = 1476
$ws.Range("D130").Value = 2187261
$ws.Range("C134").Value = 34211
$ws.Range("D134").Value = 45291943
$ws.Range("C135").Value = 14518
$ws.Range("D135").Value = 20998973
$ws.Range("C136").Value = 34618
$ws.Range("D136").Value = 50805744
$ws.Range("C137").Value = 12211
$ws.Range("D137").Value = 18142771
$ws.Range("C138").Value = 3232
$ws.Range("D138").Value = 4817875
$ws.Range("C143").Value = 11603
$ws.Range("D143").Value = 15414361
$ws.Range("C144").Value = 39222
$ws.Range("D144").Value = 56622636
$ws.Range("C145").Value = 89798
$ws.Range("D145").Value = 131454155
$ws.Range("C146").Value = 26835
$ws.Range("D146").Value = 39863611
$ws.Range("C147").Value = 7154
$ws.Range("D147").Value = 10660774
$ws.Range("C148").Value = 1750
$ws.Range("D148").Value = 2600406
$ws.Range("C151").Value = 31717
$ws.Range("D151").Value = 42621580
